$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.1535818512568701
$ws.Range("D3").Value = 0.1522119933428666
$ws.Range("D4").Value = 0.131575099332171
$ws.Range("D5").Value = 0.1316127768341236
$ws.Range("D6").Value = 0.1315393259675999
$ws.Range("D7").Value = 0.149528723519125
$ws.Range("D8").Value = 0.1499502297472438
